# Weekly fruit/vegetable price update: two new price records were added for
# "Arándano (blue)" at Vega Modelo de Temuco, inserted right after the
# existing row for 2021-12-02 (old row 10). All the subsequent records shift
# down by two rows to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 11, pushing all the
# existing data (old rows 11-68) down to rows 13-70.
$ws.Rows("11:12").Insert()

# New record #1 -> new row 11
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Vega Modelo de Temuco"
$ws.Range("C11").Value = "La Araucanía"
$ws.Range("D11").Value = 44537
$ws.Range("E11").Value = 9
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100101
$ws.Range("H11").Value = "Berries"
$ws.Range("I11").Value = 100101001
$ws.Range("J11").Value = "Arándano (blue)"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 210
$ws.Range("N11").Value = 3000
$ws.Range("O11").Value = 3000
$ws.Range("P11").Value = 3000
$ws.Range("Q11").Value = "$/kilo"
$ws.Range("R11").Value = "Región del Maule"
$ws.Range("S11").Value = 3000
$ws.Range("T11").Value = 1

# New record #2 -> new row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Vega Modelo de Temuco"
$ws.Range("C12").Value = "La Araucanía"
$ws.Range("D12").Value = 44537
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100101
$ws.Range("H12").Value = "Berries"
$ws.Range("I12").Value = 100101001
$ws.Range("J12").Value = "Arándano (blue)"
$ws.Range("K12").Value = "Sin especificar"
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 2000
$ws.Range("O12").Value = 2000
$ws.Range("P12").Value = 2000
$ws.Range("Q12").Value = "$/kilo"
$ws.Range("R12").Value = "Región del Maule"
$ws.Range("S12").Value = 2000
$ws.Range("T12").Value = 1
